$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 47
$ws.Range("H47").Value = 9999.5
$ws.Range("I47").Value = 9999.5
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 9999.5
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -9027.5
$ws.Range("N47").Value = $null

# Row 62
$ws.Range("H62").Value = 40003664
$ws.Range("I62").Value = 66669780
$ws.Range("J62").Value = 4493.6
$ws.Range("K62").Value = 66669780
$ws.Range("L62").Value = 4493.6
$ws.Range("M62").Value = -66669156
$ws.Range("N62").Value = -5741.6

# Row 65
$ws.Range("H65").Value = 40003664
$ws.Range("I65").Value = 66669780
$ws.Range("J65").Value = 4493.6
$ws.Range("K65").Value = 333348900
$ws.Range("L65").Value = 22468
$ws.Range("M65").Value = -333345780
$ws.Range("N65").Value = -28708

# Row 76
$ws.Range("H76").Value = 3365.39
$ws.Range("I76").Value = 3121.074
$ws.Range("J76").Value = 6004
$ws.Range("K76").Value = 3121.074
$ws.Range("L76").Value = 6004
$ws.Range("M76").Value = -2806.074
$ws.Range("N76").Value = -6634

# Row 79
$ws.Range("H79").Value = 3365.39
$ws.Range("I79").Value = 3121.074
$ws.Range("J79").Value = 6004
$ws.Range("K79").Value = 3121.074
$ws.Range("L79").Value = 6004
$ws.Range("M79").Value = -2029.074
$ws.Range("N79").Value = -8188

$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 2366.5813
$ws.Range("I132").Value = 1047.6364
$ws.Range("J132").Value = 3748.3333
$ws.Range("K132").Value = 3142.9092
$ws.Range("L132").Value = 11244.9999
$ws.Range("M132").Value = -612.9092000000001
$ws.Range("N132").Value = -16304.9999

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1737.5
$ws.Range("I105").Value = 1800
$ws.Range("J105").Value = 1633.3334
$ws.Range("K105").Value = 1800
$ws.Range("L105").Value = 1633.3334
$ws.Range("M105").Value = -53
$ws.Range("N105").Value = -5127.3334

# Row 134
$ws.Range("H134").Value = 3052.1
$ws.Range("I134").Value = 2102.8
$ws.Range("J134").Value = 4001.4
$ws.Range("K134").Value = 6308.400000000001
$ws.Range("L134").Value = 12004.2
$ws.Range("M134").Value = -3773.400000000001
$ws.Range("N134").Value = -17074.2

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1732.7273
$ws.Range("I16").Value = 1932.5
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1932.5
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -1645.5
$ws.Range("N16").Value = -1774

# Row 105
$ws.Range("H105").Value = 202584
$ws.Range("I105").Value = 202584
$ws.Range("K105").Value = 202584
$ws.Range("M105").Value = -200837

# Row 113
$ws.Range("H113").Value = 1732.7273
$ws.Range("I113").Value = 1932.5
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1932.5
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 237.5
$ws.Range("N113").Value = -5540

# Row 122
$ws.Range("H122").Value = 112492.89
$ws.Range("I122").Value = 167753.67
$ws.Range("J122").Value = 1971.3334
$ws.Range("K122").Value = 503261.01
$ws.Range("L122").Value = 5914.0002
$ws.Range("M122").Value = -500811.01
$ws.Range("N122").Value = -10814.0002

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 173
$ws.Range("I14").Value = 173
$ws.Range("K14").Value = 519
$ws.Range("M14").Value = -346

# Row 68
$ws.Range("H68").Value = 921.89
$ws.Range("I68").Value = 704.25
$ws.Range("J68").Value = 1092.8928
$ws.Range("K68").Value = 2112.75
$ws.Range("L68").Value = 3278.6784
$ws.Range("M68").Value = -1301.75
$ws.Range("N68").Value = -4900.678400000001

# Row 71
$ws.Range("H71").Value = 921.89
$ws.Range("I71").Value = 704.25
$ws.Range("J71").Value = 1092.8928
$ws.Range("K71").Value = 6338.25
$ws.Range("L71").Value = 9836.0352
$ws.Range("M71").Value = -2282.25
$ws.Range("N71").Value = -17948.0352

# Row 131
$ws.Range("H131").Value = 5902.9546
$ws.Range("J131").Value = 8919
$ws.Range("L131").Value = 26757
$ws.Range("N131").Value = -36837

# Row 137
$ws.Range("H137").Value = 2457.258
$ws.Range("I137").Value = 1797.3158
$ws.Range("J137").Value = 3502.1667
$ws.Range("K137").Value = 5391.9474
$ws.Range("L137").Value = 10506.5001
$ws.Range("M137").Value = -291.9474
$ws.Range("N137").Value = -20706.5001

# Row 140
$ws.Range("H140").Value = 1766.8718
$ws.Range("I140").Value = 1158.7693
$ws.Range("J140").Value = 2983.077
$ws.Range("K140").Value = 3476.3079
$ws.Range("L140").Value = 8949.231
$ws.Range("M140").Value = 1703.6921
$ws.Range("N140").Value = -19309.231

$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").Value = $null

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = $null

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = $null

# Row 80
$ws.Range("H80").Value = 4990.609
$ws.Range("I80").Value = 4932.5
$ws.Range("J80").Value = 5199.8
$ws.Range("K80").Value = 4932.5
$ws.Range("L80").Value = 5199.8
$ws.Range("M80").Value = -3934.5
$ws.Range("N80").Value = -7195.8

# Row 83
$ws.Range("H83").Value = 4990.609
$ws.Range("I83").Value = 4932.5
$ws.Range("J83").Value = 5199.8
$ws.Range("K83").Value = 24662.5
$ws.Range("L83").Value = 25999
$ws.Range("M83").Value = -19670.5
$ws.Range("N83").Value = -35983

# Row 132
$ws.Range("H132").Value = 8068.6665
$ws.Range("I132").Value = 6200
$ws.Range("J132").Value = 9003
$ws.Range("K132").Value = 18600
$ws.Range("L132").Value = 27009
$ws.Range("M132").Value = -16070
$ws.Range("N132").Value = -32069

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 58826.39
$ws.Range("I7").Value = 94433.63
$ws.Range("J7").Value = 2872.1428
$ws.Range("K7").Value = 94433.63
$ws.Range("L7").Value = 2872.1428
$ws.Range("M7").Value = -94321.63
$ws.Range("N7").Value = -3096.1428

# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = $null

# Row 119
$ws.Range("H119").Value = 42210
$ws.Range("J119").Value = 42210
$ws.Range("L119").Value = 42210
$ws.Range("N119").Value = -51886

# Row 126
$ws.Range("H126").Value = 58826.39
$ws.Range("I126").Value = 94433.63
$ws.Range("J126").Value = 2872.1428
$ws.Range("K126").Value = 283300.89
$ws.Range("L126").Value = 8616.428400000001
$ws.Range("M126").Value = -280830.89
$ws.Range("N126").Value = -13556.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 47409
$ws.Range("I126").Value = 68120
$ws.Range("J126").Value = 3028.2856
$ws.Range("K126").Value = 204360
$ws.Range("L126").Value = 9084.856800000001
$ws.Range("M126").Value = -201890
$ws.Range("N126").Value = -14024.8568
